$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2047.5454
$ws.Range("I6").Value = 2989
$ws.Range("K6").Value = 8967
$ws.Range("M6").Value = -8855

$ws.Range("H16").Value = 52000
$ws.Range("J16").Value = 52000
$ws.Range("L16").Value = 52000
$ws.Range("N16").Value = -52460

$ws.Range("H62").Value = 5318.9
$ws.Range("I62").Value = 5099.6
$ws.Range("K62").Value = 5099.6
$ws.Range("M62").Value = -4475.6

$ws.Range("H65").Value = 5318.9
$ws.Range("I65").Value = 5099.6
$ws.Range("K65").Value = 25498
$ws.Range("M65").Value = -22378

$ws.Range("H69").Value = 512.5
$ws.Range("I69").Value = 512.5
$ws.Range("K69").Value = 1537.5
$ws.Range("M69").Value = -663.5

$ws.Range("H72").Value = 512.5
$ws.Range("I72").Value = 512.5
$ws.Range("K72").Value = 4612.5
$ws.Range("M72").Value = -244.5

$ws.Range("H98").Value = 1736.091
$ws.Range("I98").Value = 1571.4375
$ws.Range("J98").Value = 2175.1667
$ws.Range("K98").Value = 1571.4375
$ws.Range("L98").Value = 2175.1667
$ws.Range("M98").Value = -73.4375
$ws.Range("N98").Value = -5171.1667

$ws.Range("H101").Value = 546.6667
$ws.Range("I101").Value = 542
$ws.Range("K101").Value = 1626
$ws.Range("M101").Value = -4

$ws.Range("H107").Value = 725.0833
$ws.Range("J107").Value = 647.3333
$ws.Range("L107").Value = 647.3333
$ws.Range("N107").Value = -4487.3333

$ws.Range("H113").Value = 4365.5454
$ws.Range("I113").Value = 4093.3333
$ws.Range("J113").Value = 4692.2
$ws.Range("K113").Value = 4093.3333
$ws.Range("L113").Value = 4692.2
$ws.Range("M113").Value = -839.3332999999998
$ws.Range("N113").Value = -11200.2

$ws.Range("H122").Value = 1736.091
$ws.Range("I122").Value = 1571.4375
$ws.Range("J122").Value = 2175.1667
$ws.Range("K122").Value = 4714.3125
$ws.Range("L122").Value = 6525.500100000001
$ws.Range("M122").Value = -2264.3125
$ws.Range("N122").Value = -11425.5001

$ws.Range("H132").Value = 3402.7334
$ws.Range("I132").Value = 3217.9285
$ws.Range("K132").Value = 9653.7855
$ws.Range("M132").Value = -7123.7855

$ws.Range("H138").Value = 1631.84
$ws.Range("I138").Value = 949
$ws.Range("J138").Value = 1761.9048
$ws.Range("K138").Value = 2847
$ws.Range("L138").Value = 5285.7144
$ws.Range("M138").Value = 2293
$ws.Range("N138").Value = -15565.7144

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1540.7142
$ws.Range("I61").Value = 1540.7142
$ws.Range("K61").Value = 1540.7142
$ws.Range("M61").Value = -1328.7142

$ws.Range("H88").Value = 1624.7916
$ws.Range("J88").Value = 2186.0625
$ws.Range("L88").Value = 2186.0625
$ws.Range("N88").Value = -2998.0625

$ws.Range("H91").Value = 1624.7916
$ws.Range("J91").Value = 2186.0625
$ws.Range("L91").Value = 2186.0625
$ws.Range("N91").Value = -4994.0625

$ws.Range("H106").Value = 10370
$ws.Range("J106").Value = 10370
$ws.Range("L106").Value = 10370
$ws.Range("N106").Value = -12894

$ws.Range("H122").Value = 2997.5
$ws.Range("I122").Value = 2997.5
$ws.Range("K122").Value = 8992.5
$ws.Range("M122").Value = -6542.5

$ws.Range("H136").Value = 1540.7142
$ws.Range("I136").Value = 1540.7142
$ws.Range("K136").Value = 4622.142599999999
$ws.Range("M136").Value = -2072.142599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5934.727
$ws.Range("I134").Value = 5932.85
$ws.Range("J134").Value = 5953.5
$ws.Range("K134").Value = 17798.55
$ws.Range("L134").Value = 17860.5
$ws.Range("M134").Value = -15263.55
$ws.Range("N134").Value = -22930.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 125
$ws.Range("I22").Value = 125
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 125
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 225
$ws.Range("N22").ClearContents()

$ws.Range("H106").Value = 22500
$ws.Range("J106").Value = 22500
$ws.Range("L106").Value = 22500
$ws.Range("N106").Value = -25024

$ws.Range("H134").Value = 7301.5
$ws.Range("J134").Value = 9400
$ws.Range("L134").Value = 28200
$ws.Range("N134").Value = -33270

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1001
$ws.Range("I92").Value = 1001
$ws.Range("K92").Value = 3003
$ws.Range("M92").Value = -1755

$ws.Range("H140").Value = 2954.2144
$ws.Range("I140").Value = 2613.3333
$ws.Range("K140").Value = 7839.999899999999
$ws.Range("M140").Value = -2659.999899999999

$ws.Range("H141").Value = 1890
$ws.Range("I141").Value = 1890
$ws.Range("K141").Value = 5670
$ws.Range("M141").Value = -490

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 950
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

$ws.Range("H70").Value = 5999.5
$ws.Range("I70").Value = 2000
$ws.Range("J70").Value = 9999
$ws.Range("K70").Value = 2000
$ws.Range("L70").Value = 9999
$ws.Range("M70").Value = -1730
$ws.Range("N70").Value = -10539

$ws.Range("H73").Value = 5999.5
$ws.Range("I73").Value = 2000
$ws.Range("J73").Value = 9999
$ws.Range("K73").Value = 2000
$ws.Range("L73").Value = 9999
$ws.Range("M73").Value = -1064
$ws.Range("N73").Value = -11871

$ws.Range("H93").Value = 90000
$ws.Range("I93").Value = 90000
$ws.Range("K93").Value = 90000
$ws.Range("M93").Value = -88128

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 3000
$ws.Range("J21").Value = 3000
$ws.Range("L21").Value = 3000
$ws.Range("N21").Value = -3348

$ws.Range("H22").Value = 1248.625
$ws.Range("I22").Value = 1012.5
$ws.Range("K22").Value = 1012.5
$ws.Range("M22").Value = -717.5

$ws.Range("H27").Value = 1248.625
$ws.Range("I27").Value = 1012.5
$ws.Range("K27").Value = 1012.5
$ws.Range("M27").Value = -905.5

$ws.Range("H46").Value = 2374.25
$ws.Range("I46").Value = 2483.9
$ws.Range("J46").Value = 2191.5
$ws.Range("K46").Value = 2483.9
$ws.Range("L46").Value = 2191.5
$ws.Range("M46").Value = -2295.9
$ws.Range("N46").Value = -2567.5

$ws.Range("H132").Value = 9632.25
$ws.Range("I132").Value = 10378.857
$ws.Range("J132").Value = 4406
$ws.Range("K132").Value = 31136.571
$ws.Range("L132").Value = 13218
$ws.Range("M132").Value = -28606.571
$ws.Range("N132").Value = -18278

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 9999999
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("N18").ClearContents()

$ws.Range("H135").Value = 56666.668
$ws.Range("J135").Value = 56666.668
$ws.Range("L135").Value = 56666.668
$ws.Range("N135").Value = -66806.66800000001
